# Energy forms and changes, graphing lines splited as individual pages
# Adds two new simulation rows (row 23: graphing-lines, row 24: energy-forms-and-changes)
# to Sheet1 of the html-simulations workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 : graphing-lines -------------------------------------------------
# Cell values are written in the same order the original author entered them so
# that any newly-introduced shared strings land at the expected indices.
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "graphing-lines"
$ws.Range("C23").Value = "YES"
$ws.Range("I23").Value = "CONVERAL TED TO INDIVIDUAL PAGES"

# --- Row 24 : energy-forms-and-changes ---------------------------------------
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "energy-forms-and-changes"
$ws.Range("C24").Value = "YES"
$ws.Range("D24").Value = "Introduction"
$ws.Range("E24").Value = "systems"
$ws.Range("I24").Value = "CONVERAL TED TO INDIVIDUAL PAGES"

# back to row 23 for the remaining OPTION columns
$ws.Range("D23").Value = "slope"
$ws.Range("E23").Value = "slope-intercept"
$ws.Range("F23").Value = "point-slope"
$ws.Range("G23").Value = "line-game"

# Column E needs to widen a little to fit "slope-intercept"
$ws.Columns.Item(5).ColumnWidth = 14

# Move / update the active selection to J23, matching the edited workbook
$ws.Range("J23").Select() | Out-Null
